$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns G and H
$ws.Range("G1").Value = "ci.lower"
$ws.Range("H1").Value = "ci.upper"

# ci.lower / ci.upper values per row (rows 12-14 intentionally left blank,
# matching the existing blank C:F cells for those rows)
$data = @(
    @(2,  0.0117359137272205,  0.418891598375897),
    @(3,  -0.0299349907384979, 0.107556059194314),
    @(4,  -0.0508819802643883, 0.0857380445883255),
    @(5,  -0.10585453306536,   0.380334055347551),
    @(6,  -0.176126317402535,  0.29677944875909),
    @(7,  -0.104724835223661,  0.376275064683585),
    @(8,  -0.178005967851821,  0.299946730244481),
    @(9,  -0.0396594531359149, 0.0958980195257914),
    @(10, 0.11333002918473,    0.154364896105377),
    @(11, 0.00116547685628543, 0.0415995272755933),
    @(15, -0.0770426803641871, 0.0635691181061893),
    @(16, -0.269526791466564,  0.222390762607537),
    @(17, -0.262910632205396,  0.219170696311658),
    @(18, -0.171476089176651,  0.101763960528777),
    @(19, -0.599893460488961,  0.356011935703642),
    @(20, -0.593558897518181,  0.352252634805069),
    @(21, 0.00116547685628543, 0.0415995272755933),
    @(22, -0.0793189062718298, 0.191796039051583),
    @(23, 0.00407731449786333, 0.145532152569401),
    @(24, -0.277489960227567,  0.67098095208015),
    @(25, 0.00588889376575046, 0.147937497159885),
    @(26, -0.276740196640265,  0.671872850279013)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 7).Value = $row[1]
    $ws.Cells.Item($r, 8).Value = $row[2]
}
